# ------------------------------------------------------------------
# Commit: "Add first 25 runes for workflow 7"
#
# Fills in the measurement columns (E..J) for the first 25 runs of
# workflow 7 (rows 302-326), adds the matching MEDIAN() rollup in O8
# and restores the workbook/sheet view state (zoom + selection) that
# Excel recorded after the edit.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# --- Workflow 7, runs 1-25: fill in the per-node measurements (E:J) ---
# Column D already carries a shared "=SUM(Ex:Jx)" formula for these rows,
# so writing E:J lets Excel recompute D automatically.
$ws.Cells.Item(302, 5).Value = 158574
$ws.Cells.Item(302, 6).Value = 160663
$ws.Cells.Item(302, 7).Value = 163370
$ws.Cells.Item(302, 8).Value = 161046
$ws.Cells.Item(302, 9).Value = 154067
$ws.Cells.Item(302, 10).Value = 150995

$ws.Cells.Item(303, 5).Value = 161658
$ws.Cells.Item(303, 6).Value = 157946
$ws.Cells.Item(303, 7).Value = 157214
$ws.Cells.Item(303, 8).Value = 157868
$ws.Cells.Item(303, 9).Value = 151516
$ws.Cells.Item(303, 10).Value = 148463

$ws.Cells.Item(304, 5).Value = 149064
$ws.Cells.Item(304, 6).Value = 153583
$ws.Cells.Item(304, 7).Value = 151478
$ws.Cells.Item(304, 8).Value = 160312
$ws.Cells.Item(304, 9).Value = 148263
$ws.Cells.Item(304, 10).Value = 156441

$ws.Cells.Item(305, 5).Value = 151316
$ws.Cells.Item(305, 6).Value = 157481
$ws.Cells.Item(305, 7).Value = 155128
$ws.Cells.Item(305, 8).Value = 162789
$ws.Cells.Item(305, 9).Value = 161158
$ws.Cells.Item(305, 10).Value = 149467

$ws.Cells.Item(306, 5).Value = 160669
$ws.Cells.Item(306, 6).Value = 155505
$ws.Cells.Item(306, 7).Value = 161425
$ws.Cells.Item(306, 8).Value = 151684
$ws.Cells.Item(306, 9).Value = 159768
$ws.Cells.Item(306, 10).Value = 156174

$ws.Cells.Item(307, 5).Value = 149744
$ws.Cells.Item(307, 6).Value = 160431
$ws.Cells.Item(307, 7).Value = 149443
$ws.Cells.Item(307, 8).Value = 154114
$ws.Cells.Item(307, 9).Value = 153019
$ws.Cells.Item(307, 10).Value = 161652

$ws.Cells.Item(308, 5).Value = 148234
$ws.Cells.Item(308, 6).Value = 148019
$ws.Cells.Item(308, 7).Value = 159474
$ws.Cells.Item(308, 8).Value = 157564
$ws.Cells.Item(308, 9).Value = 161255
$ws.Cells.Item(308, 10).Value = 157680

$ws.Cells.Item(309, 5).Value = 149370
$ws.Cells.Item(309, 6).Value = 157524
$ws.Cells.Item(309, 7).Value = 156004
$ws.Cells.Item(309, 8).Value = 152274
$ws.Cells.Item(309, 9).Value = 148584
$ws.Cells.Item(309, 10).Value = 153342

$ws.Cells.Item(310, 5).Value = 160964
$ws.Cells.Item(310, 6).Value = 155676
$ws.Cells.Item(310, 7).Value = 159099
$ws.Cells.Item(310, 8).Value = 150479
$ws.Cells.Item(310, 9).Value = 151684
$ws.Cells.Item(310, 10).Value = 158362

$ws.Cells.Item(311, 5).Value = 156352
$ws.Cells.Item(311, 6).Value = 155458
$ws.Cells.Item(311, 7).Value = 153215
$ws.Cells.Item(311, 8).Value = 158029
$ws.Cells.Item(311, 9).Value = 161011
$ws.Cells.Item(311, 10).Value = 157344

$ws.Cells.Item(312, 5).Value = 151497
$ws.Cells.Item(312, 6).Value = 156978
$ws.Cells.Item(312, 7).Value = 151806
$ws.Cells.Item(312, 8).Value = 149992
$ws.Cells.Item(312, 9).Value = 150218
$ws.Cells.Item(312, 10).Value = 154219

$ws.Cells.Item(313, 5).Value = 160242
$ws.Cells.Item(313, 6).Value = 148681
$ws.Cells.Item(313, 7).Value = 149423
$ws.Cells.Item(313, 8).Value = 156637
$ws.Cells.Item(313, 9).Value = 154871
$ws.Cells.Item(313, 10).Value = 160936

$ws.Cells.Item(314, 5).Value = 157736
$ws.Cells.Item(314, 6).Value = 149070
$ws.Cells.Item(314, 7).Value = 154291
$ws.Cells.Item(314, 8).Value = 148256
$ws.Cells.Item(314, 9).Value = 150892
$ws.Cells.Item(314, 10).Value = 158557

$ws.Cells.Item(315, 5).Value = 150271
$ws.Cells.Item(315, 6).Value = 159752
$ws.Cells.Item(315, 7).Value = 154253
$ws.Cells.Item(315, 8).Value = 151376
$ws.Cells.Item(315, 9).Value = 157554
$ws.Cells.Item(315, 10).Value = 149511

$ws.Cells.Item(316, 5).Value = 148357
$ws.Cells.Item(316, 6).Value = 152008
$ws.Cells.Item(316, 7).Value = 152935
$ws.Cells.Item(316, 8).Value = 160283
$ws.Cells.Item(316, 9).Value = 160979
$ws.Cells.Item(316, 10).Value = 148073

$ws.Cells.Item(317, 5).Value = 151033
$ws.Cells.Item(317, 6).Value = 157933
$ws.Cells.Item(317, 7).Value = 156175
$ws.Cells.Item(317, 8).Value = 156164
$ws.Cells.Item(317, 9).Value = 152760
$ws.Cells.Item(317, 10).Value = 160456

$ws.Cells.Item(318, 5).Value = 151476
$ws.Cells.Item(318, 6).Value = 154393
$ws.Cells.Item(318, 7).Value = 149645
$ws.Cells.Item(318, 8).Value = 161448
$ws.Cells.Item(318, 9).Value = 154036
$ws.Cells.Item(318, 10).Value = 148314

$ws.Cells.Item(319, 5).Value = 157698
$ws.Cells.Item(319, 6).Value = 153526
$ws.Cells.Item(319, 7).Value = 161364
$ws.Cells.Item(319, 8).Value = 149780
$ws.Cells.Item(319, 9).Value = 157933
$ws.Cells.Item(319, 10).Value = 153742

$ws.Cells.Item(320, 5).Value = 149531
$ws.Cells.Item(320, 6).Value = 148230
$ws.Cells.Item(320, 7).Value = 161324
$ws.Cells.Item(320, 8).Value = 152685
$ws.Cells.Item(320, 9).Value = 160522
$ws.Cells.Item(320, 10).Value = 160445

$ws.Cells.Item(321, 5).Value = 156415
$ws.Cells.Item(321, 6).Value = 154196
$ws.Cells.Item(321, 7).Value = 148861
$ws.Cells.Item(321, 8).Value = 160633
$ws.Cells.Item(321, 9).Value = 152437
$ws.Cells.Item(321, 10).Value = 150039

$ws.Cells.Item(322, 5).Value = 148442
$ws.Cells.Item(322, 6).Value = 155009
$ws.Cells.Item(322, 7).Value = 161566
$ws.Cells.Item(322, 8).Value = 154771
$ws.Cells.Item(322, 9).Value = 153698
$ws.Cells.Item(322, 10).Value = 156600

$ws.Cells.Item(323, 5).Value = 154541
$ws.Cells.Item(323, 6).Value = 156599
$ws.Cells.Item(323, 7).Value = 158352
$ws.Cells.Item(323, 8).Value = 156875
$ws.Cells.Item(323, 9).Value = 154690
$ws.Cells.Item(323, 10).Value = 159134

$ws.Cells.Item(324, 5).Value = 151011
$ws.Cells.Item(324, 6).Value = 160866
$ws.Cells.Item(324, 7).Value = 154397
$ws.Cells.Item(324, 8).Value = 153190
$ws.Cells.Item(324, 9).Value = 160095
$ws.Cells.Item(324, 10).Value = 148599

$ws.Cells.Item(325, 5).Value = 156962
$ws.Cells.Item(325, 6).Value = 160798
$ws.Cells.Item(325, 7).Value = 161421
$ws.Cells.Item(325, 8).Value = 159997
$ws.Cells.Item(325, 9).Value = 157457
$ws.Cells.Item(325, 10).Value = 158421

$ws.Cells.Item(326, 5).Value = 153998
$ws.Cells.Item(326, 6).Value = 148838
$ws.Cells.Item(326, 7).Value = 159790
$ws.Cells.Item(326, 8).Value = 152521
$ws.Cells.Item(326, 9).Value = 159086
$ws.Cells.Item(326, 10).Value = 156772

# --- New rollup cell: median of the (now non-empty) D302:D351 block ---
$ws.Range("O8").Formula = "=MEDIAN(D302:D351)"

# --- Sheet view: scroll/zoom/selection as left by the editor ---
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 313
$win.ScrollColumn = 1
$ws.Range("N5").Select()
$win.Zoom = 85

# --- Workbook window geometry as left by the editor ---
$win.Left = -105
$win.Top = 0
$win.Width = 19410
$win.Height = 15585

Write-Host "Workflow 7 runs 1-25 populated; O8 rollup added."
